$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 (the existing row 3 "a,a,a,a,a" shifts down to row 4)
$ws.Rows.Item(3).Insert()

# Update row 2: B2 becomes empty text, C2 becomes a text CPF, D2/E2 updated
$ws.Cells.Item(2, 2).Value = ""
$ws.Cells.Item(2, 3).Value = "123.456.789-00"
$ws.Cells.Item(2, 4).Value = "artur"
$ws.Cells.Item(2, 5).Value = 123

# Fill the newly inserted row 3: A3/B3/C3 empty text, D3 "artur", E3 = 1
$ws.Cells.Item(3, 1).Value = ""
$ws.Cells.Item(3, 2).Value = ""
$ws.Cells.Item(3, 3).Value = ""
$ws.Cells.Item(3, 4).Value = "artur"
$ws.Cells.Item(3, 5).Value = 1
